$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the numeric filler values that used to occupy J2:S2
$ws.Range("J2:S2").Clear()

# Row 3: add the two new xpath-description columns
$ws.Range("M3").Value = "Ingresa un destino"
$ws.Range("N3").Value = "Ingresa una fecha de partida"

# Row 4: new "Ofertas" test case data
$ws.Range("C4").Value = "agosto 2023"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("D4").Value = "Orlando"
$ws.Range("E4").Value = "Vuelo a Orlando"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = "Despegar - Checkout"

# Row 5: another "Ofertas" test case data
$ws.Range("B5").Value = "Ofertas"
$ws.Range("C5").Value = "agosto 2024"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("D5").Value = "Londres"

# Update the active selection to reflect where the user ended up
$ws.Range("C5").Select()
